$wb = $excel.ActiveWorkbook

# "funding" sheet (sheet7.xml): update the funder row.
$fundingWs = $wb.Worksheets.Item("funding")

$fundingWs.Range("A2").Value = "California Department of Water Resources "
$fundingWs.Range("B2").Value = "DWR"
$fundingWs.Range("C2").Value = "NA"
$fundingWs.Range("D2").Value = "NA"
$fundingWs.Range("F2").ClearContents()

# The "funding" sheet becomes the active/selected tab, with E6 selected -
# moving tab-selection away from "keyword_set".
$fundingWs.Activate() | Out-Null
$fundingWs.Range("E6").Select() | Out-Null
